$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (row 19) so the table shrinks from 18 to 17 data rows
$ws.Rows(19).Delete()

# New data for the table (player, position(s), team), in the final row order
$data = @(
  @("Dyson Daniels","PG,SG,SF","Atlanta Hawks"),
  @("Malik Beasley","SG,SF","Detroit Pistons"),
  @("Josh Hart","SG,SF,PF","New York Knicks"),
  @("Bilal Coulibaly","SG,SF","Washington Wizards"),
  @("Michael Porter Jr.","SF,PF","Denver Nuggets"),
  @("Jaden McDaniels","SF,PF","Minnesota Timberwolves"),
  @("De'Andre Hunter","SF,PF","Atlanta Hawks"),
  @("Andrew Wiggins","SF,PF","Golden State Warriors"),
  @("Domantas Sabonis","C","Sacramento Kings"),
  @("Victor Wembanyama","C","San Antonio Spurs"),
  @("Alperen Sengün","C","Houston Rockets"),
  @("Kelly Oubre Jr.","SG,SF","Philadelphia 76ers"),
  @("Donovan Mitchell","PG,SG","Cleveland Cavaliers"),
  @("Kristaps Porzingis","PF,C","Boston Celtics"),
  @("Kel'el Ware","PF,C","Miami Heat"),
  @("Cam Thomas","SG,SF","Brooklyn Nets"),
  @("Donte DiVincenzo","PG,SG,SF","Minnesota Timberwolves")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
  $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
